# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) timestamps for rows 33-34
# on the "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D33:D34").Value = "2016-03-07 02:52:35"
$wsZhCn.Range("G33:G34").Value = "2016-03-07 02:53:22"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D33:D34").Value = "2016-03-07 02:52:45"
$wsDeDe.Range("G33:G34").Value = "2016-03-07 02:53:42"
